# Add a new "description" column header to the "params" sheet.
#
# This mirrors the commit "added checks for no description, invalid header":
# a new `description` column (M) is appended after the existing header row
# on the `params` sheet, which also registers a new shared string and grows
# the sheet's used range / row spans from A1:L3 to A1:M3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")
$ws.Activate()

# Add the new header cell with the "description" text (creates a new shared
# string entry and extends the sheet dimension/row spans automatically).
$ws.Range("M1").Value = "description"

# Reflect the cursor/selection moving to the newly added header cell, and
# scroll the view right so the new column is visible, matching the author's
# editing session.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M1").Select()
